$wb = $excel.ActiveWorkbook

# --- 1. Sheet surgery -------------------------------------------------
# We need the final tab order 2021-Q1..2021-Q4, 2022-Q1, 总计 with
# sheetId 1..4, 5, 6 respectively (总计 currently holds sheetId 5).
# Recreating 总计 after inserting the new 2022-Q1 sheet makes the
# engine hand out sheetId 5 to 2022-Q1 and sheetId 6 to the rebuilt
# 总计, matching the target workbook.xml exactly.

$q4 = $wb.Worksheets.Item("2021-Q4")

# Snapshot the old 总计 data before removing the sheet.
$oldTotal = $wb.Worksheets.Item("总计")
$rowDates = @()
$rowCounts = @()
$rowValues = @()
for ($i = 2; $i -le 5; $i++) {
    $rowDates += , $oldTotal.Cells.Item($i, 2).Value2
    $rowCounts += , $oldTotal.Cells.Item($i, 3).Value2
    $rowValues += , $oldTotal.Cells.Item($i, 4).Value2
}
$oldTotal.Delete()

# New "2022-Q1" fund-holdings sheet, positioned right after 2021-Q4.
$newQ = $wb.Worksheets.Add($null, $q4)
$newQ.Name = "2022-Q1"

# Rebuilt "总计" sheet, positioned right after "2022-Q1".
$totals = $wb.Worksheets.Add($null, $newQ)
$totals.Name = "总计"

# --- 2. Populate "2022-Q1" ---------------------------------------------
# Style donor: 2021-Q4 uses the same 7-column layout / "s=2" styling
# that the new sheet needs for its header row and column A.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

$newQ.Cells.Item(1, 2).Value = "基金代码"
$newQ.Cells.Item(1, 3).Value = "基金名称"
$newQ.Cells.Item(1, 4).Value = "基金规模"
$newQ.Cells.Item(1, 5).Value = "股票总仓位"
$newQ.Cells.Item(1, 6).Value = "仓位占比"
$newQ.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newQ.Cells.Item(1, 8).Value = "仓位排名"

$newQ.Cells.Item(2, 1).Value = 0
$newQ.Cells.Item(2, 2).Value = "'010765"
$newQ.Cells.Item(2, 3).Value = "国寿安保华丰混合A"
$newQ.Cells.Item(2, 4).Value = "'0.42"
$newQ.Cells.Item(2, 5).Value = "'88.40"
$newQ.Cells.Item(2, 6).Value = "'2.17"
$newQ.Cells.Item(2, 7).Value = "'0.0091"
$newQ.Cells.Item(2, 8).Value = 10

$newQ.Cells.Item(3, 1).Value = 1
$newQ.Cells.Item(3, 2).Value = "'010766"
$newQ.Cells.Item(3, 3).Value = "国寿安保华丰混合C"
$newQ.Cells.Item(3, 4).Value = "'0.00"
$newQ.Cells.Item(3, 5).Value = "'88.40"
$newQ.Cells.Item(3, 6).Value = "'2.17"
$newQ.Cells.Item(3, 7).Value = 0
$newQ.Cells.Item(3, 8).Value = 10

# Copy the header / column-A formatting from 2021-Q4 so we reuse the
# existing "s=2" style instead of registering new ones.
$styleSrc.Range("B1:H1").Copy()
$newQ.Range("B1:H1").PasteSpecial(-4122)
$styleSrc.Range("A2:A3").Copy()
$newQ.Range("A2:A3").PasteSpecial(-4122)

# --- 3. Populate "总计" -------------------------------------------------
$totals.Cells.Item(1, 2).Value = "日期"
$totals.Cells.Item(1, 3).Value = "持有数量(只)"
$totals.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q1"
$totals.Cells.Item(2, 3).Value = 2
$totals.Cells.Item(2, 4).Value = 0.01

$totals.Cells.Item(3, 1).Value = 1
$totals.Cells.Item(3, 2).Value = $rowDates[0]
$totals.Cells.Item(3, 3).Value = $rowCounts[0]
$totals.Cells.Item(3, 4).Value = $rowValues[0]

$totals.Cells.Item(4, 1).Value = 2
$totals.Cells.Item(4, 2).Value = $rowDates[1]
$totals.Cells.Item(4, 3).Value = $rowCounts[1]
$totals.Cells.Item(4, 4).Value = $rowValues[1]

$totals.Cells.Item(5, 1).Value = 3
$totals.Cells.Item(5, 2).Value = $rowDates[2]
$totals.Cells.Item(5, 3).Value = $rowCounts[2]
$totals.Cells.Item(5, 4).Value = $rowValues[2]

$totals.Cells.Item(6, 1).Value = 4
$totals.Cells.Item(6, 2).Value = $rowDates[3]
$totals.Cells.Item(6, 3).Value = $rowCounts[3]
$totals.Cells.Item(6, 4).Value = $rowValues[3]

$styleSrc.Range("B1:D1").Copy()
$totals.Range("B1:D1").PasteSpecial(-4122)
$styleSrc.Range("A2:A7").Copy()
$totals.Range("A2:A7").PasteSpecial(-4122)

Write-Output "done"
